$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (bold, centered, bordered) from an
# existing header cell so the new header cells match the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every data row (2-57): Wins=89, Losses=73, Ties=0
for ($row = 2; $row -le 57; $row++) {
    $ws.Cells.Item($row, 30).Value = 89
    $ws.Cells.Item($row, 31).Value = 73
    $ws.Cells.Item($row, 32).Value = 0
}
